# Fix typo in shared string used by cell E1 ("Дата предоставление кредита"
# -> "Дата предоставления кредита")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("Дата предоставление кредита", "Дата предоставления кредита") | Out-Null

# Apply the (new) default-looking style to the previously-blank cells:
# F1:F3 and the whole A4:F11 block move from style 2 to a brand new style 3
# (same number format/protection, but a fresh Calibri 11 "automatic colour"
# font instead of font 4 used by style 2).
$targetRange = $ws.Range("F1:F3")
$targetRange.Font.Name = "Calibri"
$targetRange.Font.Size = 11
$targetRange.Font.ThemeColor = 1

$targetRange2 = $ws.Range("A4:F11")
$targetRange2.Font.Name = "Calibri"
$targetRange2.Font.Size = 11
$targetRange2.Font.ThemeColor = 1

# Move the active selection from O14 to E1
$ws.Range("E1").Select() | Out-Null
